$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.149.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "'1.657.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "'217.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "'0.5192"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.71%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").Value = "'0.2657"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "'0.06279"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("D10").Value = "'20.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.31%  "
$ws.Range("D11").Value = "'0.07721"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("D12").Value = "'1.679.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("D13").Value = "'4.409"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("D14").Value = "'1.886.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.72%  "
$ws.Range("D15").Value = "'0.5430"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.18%  "
$ws.Range("D16").Value = "'0.0₅8158"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.04%  "
$ws.Range("D17").Value = "'64.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.46%  "
$ws.Range("D18").Value = "'26.173.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("D19").Value = "'1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").Value = "'4.640"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.45%  "
$ws.Range("D21").Value = "'191.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("D22").Value = "'10.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.05%  "
$ws.Range("D23").Value = "'6.050"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.25%  "
$ws.Range("D24").Value = "'1.008"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").Value = "'138.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.13%  "
$ws.Range("D26").Value = "'0.1232"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.18%  "
$ws.Range("D27").Value = "'7.170"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.19%  "
$ws.Range("D28").Value = "'16.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("D29").Value = "'1.409"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("D30").Value = "'0.05960"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.18%  "
$ws.Range("D31").Value = "'1.277"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").Value = "'3.556"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.38%  "
$ws.Range("D33").Value = "'3.301"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.91%  "
$ws.Range("D34").Value = "'1.630"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.26%  "
$ws.Range("D35").Value = "'0.9740"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.94%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.411"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'2.777"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").Value = "'0.5867"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("D39").Value = "'0.01578"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.31%  "
$ws.Range("D40").Value = "'5.913"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("D41").Value = "'0.8590"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").Value = "'1.004"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").Value = "'1.032.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.20%  "
$ws.Range("D44").Value = "'99.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("D45").Value = "'1.800.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.11%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.0₈108"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'56.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("D48").Value = "'1.005"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").Value = "'8.026"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").Value = "'0.05180"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").Value = "'0.4228"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
